$d = $word.ActiveDocument

# Locate the paragraph holding "Baz changes" (it also holds the _GoBack
# bookmark, inline, between "Baz chan" and "ges").
$hit = $d.Content
$hit.Find.Execute("Baz changes", $false, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0)

$targetIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Start -le $hit.Start -and $hit.Start -lt $candidate.Range.End) {
        $targetIndex = $i
        break
    }
}

# Paragraph 5 ("Baz changes" with inline _GoBack bookmark) is replaced by the
# new GitHub tutorial paragraph, and the _GoBack bookmark moves into its own
# (now-empty) paragraph immediately after it -- which swallows the following
# blank paragraph so the total paragraph count stays the same.
$p5 = $d.Paragraphs.Item($targetIndex)
$p6 = $d.Paragraphs.Item($targetIndex + 1)
$r = $d.Range($p5.Range.Start, $p6.Range.End)

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>GitHub</w:t></w:r><w:r><w:t> is a code hosting platform for </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">version </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>control</w:t></w:r><w:r><w:t>and</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> collaboration. It lets you and others work together on projects from anywhere. This tutorial teaches </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>you</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>GitHub</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t> essentials like repositories, branches, commits, and Pull Requests</w:t></w:r></w:p><w:p><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$r.InsertXML($xml)
